$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "华胜天成"
$ws.Range("A3").Value = "拓普集团"
$ws.Range("C3").Value = "万通发展"
$ws.Range("B4").Value = "巨轮智能"
$ws.Range("C4").Value = "浙江荣泰"
$ws.Range("A5").Value = "岩山科技"
$ws.Range("B5").Value = "山子高科"
$ws.Range("C5").Value = "山子高科"
$ws.Range("A6").Value = "山子高科"
$ws.Range("B6").Value = "岩山科技"
$ws.Range("C6").Value = "利欧股份"
$ws.Range("A7").Value = "华胜天成"
$ws.Range("C7").Value = "巨轮智能"
$ws.Range("A8").Value = "巨轮智能"
$ws.Range("B8").Value = "吉视传媒"
$ws.Range("A9").Value = "浙江荣泰"
$ws.Range("C9").Value = "三维通信"
$ws.Range("A10").Value = "国轩高科"
$ws.Range("B10").Value = "数据港"
$ws.Range("C10").Value = "卧龙电驱"
$ws.Range("B11").Value = "长城军工"
$ws.Range("C11").Value = "吉视传媒"
$ws.Range("A12").Value = "太平洋"
$ws.Range("B12").Value = "秦川机床"
$ws.Range("C12").Value = "拓普集团"
$ws.Range("A13").Value = "三维通信"
$ws.Range("B13").Value = "拓普集团"
$ws.Range("C13").Value = "科森科技"
$ws.Range("A14").Value = "长城军工"
$ws.Range("B14").Value = "国轩高科"
$ws.Range("C14").Value = "长城军工"
$ws.Range("A15").Value = "万通发展"
$ws.Range("B15").Value = "供销大集"
$ws.Range("C15").Value = "天普股份"
$ws.Range("A16").Value = "天普股份"
$ws.Range("B16").Value = "三维通信"
$ws.Range("C16").Value = "领益智造"
$ws.Range("A17").Value = "新易盛"
$ws.Range("B17").Value = "上纬新材"
$ws.Range("C17").Value = "上海电力"
$ws.Range("A18").Value = "卧龙电驱"
$ws.Range("B18").Value = "卧龙电驱"
$ws.Range("C18").Value = "供销大集"
$ws.Range("A19").Value = "上纬新材"
$ws.Range("B19").Value = "新易盛"
$ws.Range("C19").Value = "恒宝股份"
$ws.Range("A20").Value = "秦川机床"
$ws.Range("B20").Value = "东方财富"
$ws.Range("C20").Value = "寒武纪"
$ws.Range("A21").Value = "供销大集"
$ws.Range("B21").Value = "科森科技"
$ws.Range("C21").Value = "兆新股份"
